# Edit script for LOQ4049.docx
# Splits several concatenated multi-sentence paragraphs into separate
# lines joined by manual line breaks (<w:br/>), matching the upstream
# "Build site" re-render of the source Markdown -> DOCX pipeline.

$d = $word.ActiveDocument

# --- Edit 1: "Objetivos" (Portuguese) -------------------------------
$p1 = $d.Paragraphs(6)
$rng1 = $p1.Range
$xml1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">1) Formativos: Propiciar ao educando as condições básicas e necessárias para a sua formação profissional. </w:t><w:br/><w:t>2) Informativos: fornecer ao educando os conceitos básicos para o entendimento, assessoramento e acompanhamento de Projetos na Indústria Química seguindo metodologia especifica.</w:t><w:br/><w:t>3) Automatizantes: desenvolver no educando o raciocínio analítico, obedecendo metodologia sistemática aplicada em projetos.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($xml1)

# --- Edit 2: "Objetivos" (English, italic) --------------------------
$p2 = $d.Paragraphs(7)
$rng2 = $p2.Range
$xml2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>1) Formative: Provide students with basic and necessary conditions for them professional qualification.</w:t><w:br/><w:t>2) Informational: Give to students basic concepts for them uptake, assistance and monitoring in Chemical Industry Design following a specifies methodology.</w:t><w:br/><w:t>3) Automated: Develop in students the reasoning analytical, following the systematic methodology applied in projects.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng2.InsertXML($xml2)

# --- Edit 3: "Programa" (Portuguese, 13 numbered items) -------------
$p3 = $d.Paragraphs(14)
$rng3 = $p3.Range
$xml3 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>1. Introdução: Conceitos de Gestão de Projetos</w:t><w:br/><w:t>2. Aspectos da Implantação de Projetos: Etapas Fundamentais e Formas Parciais</w:t><w:br/><w:t>3. O Ciclo de Vida do Projeto</w:t><w:br/><w:t>4. Aspectos da Viabilidade de Projetos: Receitas, Custos, Ponto de Nivelamento, Estimativas</w:t><w:br/><w:t>5. Guia PMBOK: Principais Áreas de Conhecimento</w:t><w:br/><w:t>6. Plano de Projeto</w:t><w:br/><w:t>7. O Gerente de Projeto e as Interfaces: Equipes de Projeto</w:t><w:br/><w:t>8. Legalização de Projetos: Aspectos sobre o Licenciamento</w:t><w:br/><w:t>9. Gestão de Riscos - Técnicas de Análise de Riscos</w:t><w:br/><w:t>10. Aspectos sobre Auditorias e Auditorias Integradas</w:t><w:br/><w:t>11. Provas e/ou apresentações de Trabalhos</w:t><w:br/><w:t>12. Sistema de Gestão Integrada (SGI): Qualidade; Meio Ambiente; e Saúde e Segurança.</w:t><w:br/><w:t>13. Legislação brasileira aplicada ao SGI</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng3.InsertXML($xml3)

# --- Edit 4: "Avaliacao" / "Metodo" run (split out "Justificativa") -
$foundMetodo = $d.Content.Find.Execute("em equipes.Justificativa", $true, $false, $false, $false, $false, $true, 1, $false, "em equipes.^lJustificativa", 2)

# --- Edit 5: "Bibliografia" ------------------------------------------
$p5 = $d.Paragraphs(19)
$rng5 = $p5.Range
$xml5 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>1) Chemical Engineering Plant Design (Vilbrandt e Bryden) 2) Project Engineering of Process Plants Rose e Barrow (2ª impressão - 1968) 3) Elaboração e Análise de Projetos Simonsen, M.H. e H. Flanger 4) Implantação de Indústrias. Valle, E.C. Livros Técnicos e Científicos Editores S/A, Rio de Janeiro. 5) Gestão Integrada: Qualidade, Meio Ambiente, Prevenção. Antecipação de riscos e outras ferramentas para implantação. CHAVES. F. J. M., 1ª Ed. 2022.</w:t><w:br/><w:t>6) PMBOK Guide (6ª Ed., 2017)</w:t><w:br/><w:t>7) Normas ABNT NBR ISO: 9001; 14001; 45001. 8) IBGR – Instituto Brasileiro de Gerenciamento de Riscos, 2000.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng5.InsertXML($xml5)

Write-Output "done"
